# Apply the evaluation update to the "EVALUACION2" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EVALUACION2")

# Row 15 ("3. Genera evidencias ...") moves from "No logrado" to "Logrado".
$ws.Range("C15").Value = "Logrado"

# Row 19 ("7.- Generan evidencias claras ...") moves from "No logrado" to
# "Completamente logrado".
$ws.Range("C19").Value = "Completamente logrado"

# Update the sheet's active selection (was F24, now E24).
[void]$ws.Activate()
[void]$ws.Range("E24").Select()
